$d = $word.ActiveDocument

# Rebuild the document body: drop the "Hello World" paragraph (and its
# Arial Black formatting on both the run and the paragraph mark), keeping
# only the _GoBack bookmark that used to trail it.
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
  '</w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xmlFrag) | Out-Null
